$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the F1 header from "Area_Name" to "Region"
$ws.Range("F1").Value = "Region"

# Insert a new row at row 18 (pushes Purdue etc. down by one) and
# populate it with the University of Illinois - Chicago chapter entry.
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = "Illinois"
$ws.Range("B18").Value = "1200 W Harrison St, Chicago, IL 60607"
$ws.Range("C18").Value = "University of Illinois - Chicago"
$ws.Range("D18").Value = "1200 W Harrison St, Chicago, IL 60607"
$ws.Range("E18").Value = "UIC"
$ws.Range("F18").Value = "Midwest"
$ws.Range("G18").Value = "None"
$ws.Range("H18").Value = "None"
$ws.Range("I18").Value = " "
